$d = $word.ActiveDocument

# Add the three new character styles referenced by the new runs.
$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.NameAscii = "Calibri"
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.NameAscii = "Calibri"
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.NameAscii = "Calibri"
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# Apply GaNStyle to every occurrence of the campaign-period sentence.
$rng = $d.Content
$found = $rng.Find.Execute("Perioadele campaniei din 2022 pentru Constelația Taurului: 16-25 ianuarie", $true)
while ($found) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $found = $rng.Find.Execute("Perioadele campaniei din 2022 pentru Constelația Taurului: 16-25 ianuarie", $true)
}

# Apply GaNParagraph to the introductory description paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("Prin această activitate participați în cadrul unei campanii globale de observare și consemnare a celor mai slabe stele vizibile ca metodă de măsurare a poluării luminoase dintr-un anumit loc. Localizând și observând  Constelația Taurului pe cerul nopții și comparând-o cu diagramele stelare, oamenii din întreaga lume vor putea afla în ce măsură iluminatul nocturn din comunitatea lor contribuie la poluarea luminoasă. Contribuțiile dumneavoastră la baza de date online vor facilita o documentare globală privind cerul nocturn observabil.", $true)
if ($found) {
    $rng.Style = "GaNParagraph"
}

# Apply GaNLinks to the credit/link line.
$rng = $d.Content
$found = $rng.Find.Execute("de Jan Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $true)
if ($found) {
    $rng.Style = "GaNLinks"
}
